$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Step 1: insert two new blank columns before column D. This shifts the
# existing columns D:K (values + per-cell styles/number formats) to F:M,
# matching two newly reported quarters being added to the left of the table. ---
$ws.Range("D1:E1").EntireColumn.Insert()

# --- Step 2: the newly inserted D:E columns have no number format yet; bulk-copy
# the numeric-data style from column F (rows 5-102) onto D:E, then fix up the three
# "Period Ending" date rows (7, 38, 80) which need the date style instead. ---
$ws.Range("F5:F102").Copy() | Out-Null
$ws.Range("D5:E102").PasteSpecial(-4122) | Out-Null
$ws.Range("G7").Copy() | Out-Null
$ws.Range("D7:E7").PasteSpecial(-4122) | Out-Null
$ws.Range("G38").Copy() | Out-Null
$ws.Range("D38:E38").PasteSpecial(-4122) | Out-Null
$ws.Range("G80").Copy() | Out-Null
$ws.Range("D80:E80").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# --- Step 3: rows that have no data in D:K (only a label in A/B) got stray blank
# D/E cells from the column insert above; clear them so the row stays label-only. ---
$ws.Range("D5:E6").Clear()
$ws.Range("D37:E37").Clear()
$ws.Range("D79:E79").Clear()

# --- Step 4: write the final figures (2 new quarters + restated prior-quarter
# values) for every data row across D:M. ---
function Set-RowValues([int]$row, [object[]]$vals) {
    $arr = New-Object 'object[,]' 1,10
    for ($i = 0; $i -lt 10; $i++) { $arr[0,$i] = $vals[$i] }
    $ws.Range("D$row`:M$row").Value2 = $arr
}

Set-RowValues 7 @(43465,43373,43281,43190,43100,43008,42916,42825,42735,42643)
Set-RowValues 8 @(20711900,21026500,19579700,24359100,20770000,20681700,18881600,23443000,19208300,19720700)
Set-RowValues 9 @(15083700,15305500,14341000,17804700,15290000,15148400,13830300,17433900,14170500,14505800)
Set-RowValues 10 @(5628200,5721000,5238700,6554300,5479900,5533300,5051300,6009100,5037800,5214900)
Set-RowValues 12 @("NA","NA","NA","NA","NA","NA","NA","NA","NA","NA")
Set-RowValues 13 @(0,0,0,0,0,0,0,0,0,0)
Set-RowValues 14 @(0,0,0,0,0,0,0,0,0,0)
Set-RowValues 15 @(0,0,0,0,0,0,0,0,0,0)
Set-RowValues 17 @(18994700,19250600,18240500,22188700,19221300,19132200,17689700,21707800,18186600,18780900)
Set-RowValues 18 @(1717100,1775800,1339300,2170300,1548600,1549500,1191800,1735100,1021800,939800)
Set-RowValues 20 @(-2082400,-171000,333600,-583700,82700,-132600,132400,-702400,286900,28100)
Set-RowValues 21 @(461200,2411300,2488000,2422700,2464400,2253700,2112600,1868900,2102600,1987900)
Set-RowValues 22 @(52000,45800,41400,45400,48500,47200,44600,38100,42000,41800)
Set-RowValues 23 @(-417300,1559000,1631400,1541300,1582800,1369600,1279600,994600,1266600,926100)
Set-RowValues 24 @(404900,396300,442200,137800,394200,340500,318200,339400,295300,179000)
Set-RowValues 25 @(0,0,0,0,0,0,0,0,0,0)
Set-RowValues 26 @(-822200,1162700,1189200,1403500,1188600,1029200,961500,655200,971300,747200)
Set-RowValues 27 @(-994100,882600,898300,1089500,873200,784200,679200,399800,697700,509700)
Set-RowValues 28 @(0,0,0,0,0,0,0,0,0,0)
Set-RowValues 29 @(-3600,-89300,53100,-145700,12400,-10900,-600,-45100,-7900,-3600)
Set-RowValues 30 @(0,0,0,0,0,0,0,0,0,0)
Set-RowValues 31 @(0,0,0,0,0,0,0,0,0,0)
Set-RowValues 32 @(2082400,171000,-333600,583700,-82700,132600,-132400,702400,-286900,-28100)
Set-RowValues 33 @(-997700,793300,951400,943800,885600,773300,678600,354800,689800,506000)
Set-RowValues 34 @(0,0,0,0,0,0,0,0,0,0)
Set-RowValues 35 @(-997700,793300,951400,943800,885600,773300,678600,354800,689800,506000)
Set-RowValues 38 @(43465,43373,43281,43190,43100,43008,42916,42825,42735,42643)
Set-RowValues 41 @(6407100,7463400,7362200,6309600,7334400,6944300,7349700,6787700,6771500,7553100)
Set-RowValues 42 @(2909200,2717100,3434500,3374800,2939900,2708500,3501600,"NA","NA","NA")
Set-RowValues 43 @(20002500,20232700,19370800,22612800,20017900,21216900,19871600,21957800,19728100,25974400)
Set-RowValues 44 @(14341300,13448300,13069600,12432100,14097600,12782800,12369500,10873800,12655400,11737600)
Set-RowValues 45 @(2261600,2108400,2145000,1842900,2258300,2187200,2212300,4753800,4708100,4654500)
Set-RowValues 46 @(45921600,45969800,45382000,46572300,46648100,45839800,45304800,44373100,43863100,49919500)
Set-RowValues 47 @(11749000,12390700,12566800,13196900,13425300,13164700,12783700,12515200,13097300,16937400)
Set-RowValues 48 @(17610700,19465500,19250700,19208400,19232200,18881700,18246800,17725900,17662900,20069600)
Set-RowValues 49 @(8952000,9700100,9568800,9531500,10139000,10111900,8483100,8153300,8193600,8309100)
Set-RowValues 50 @(0,0,0,0,0,0,0,0,0,0)
Set-RowValues 51 @(0,0,0,0,0,0,0,0,0,0)
Set-RowValues 52 @(3395900,3417900,3040400,2854500,2728600,2784700,2792500,2951400,2691400,9247200)
Set-RowValues 53 @(0,0,0,0,0,0,0,0,0,0)
Set-RowValues 54 @(87629200,90944000,89808700,91363700,92173200,90782700,87610900,85718900,85508400,104482800)
Set-RowValues 57 @(12243700,12396100,12648200,13894300,13397100,12889700,12346100,12437800,11908500,11654900)
Set-RowValues 58 @(5199600,4387200,4206600,2157200,3742700,3865500,4202800,3429100,3937900,12512300)
Set-RowValues 59 @(16408300,17095200,16671100,18258800,17693700,17628000,16719600,17137200,16265100,16240700)
Set-RowValues 60 @(33851600,33878500,33526000,34310400,34833500,34383200,33268500,33004000,32111400,40407900)
Set-RowValues 61 @(6776700,6473500,6718300,7337400,7529400,7629600,7261000,7007400,7631500,19851800)
Set-RowValues 62 @(8143600,8401900,8528600,8930400,9587400,9519400,9427600,9367200,9692000,9942400)
Set-RowValues 63 @(0,0,0,0,0,0,0,0,0,0)
Set-RowValues 64 @(0,0,0,0,0,0,0,0,0,0)
Set-RowValues 65 @(0,0,0,0,0,0,0,0,0,0)
Set-RowValues 66 @(59223000,59762700,59459600,61730400,63059800,62347900,60372500,59400900,59920500,81205500)
Set-RowValues 68 @(0,0,0,0,0,0,0,0,0,0)
Set-RowValues 69 @(0,0,0,0,0,0,0,0,0,0)
Set-RowValues 70 @(0,0,0,0,0,0,0,0,0,0)
Set-RowValues 71 @(0,0,0,0,0,0,0,0,0,0)
Set-RowValues 72 @(19295300,20630600,19838600,19032800,18040300,17443200,16664700,15909000,15530900,15048500)
Set-RowValues 73 @(0,0,0,0,0,0,0,0,0,0)
Set-RowValues 74 @(0,0,0,0,0,0,0,0,0,0)
Set-RowValues 75 @(0,0,0,0,0,0,0,0,0,0)
Set-RowValues 76 @(28406200,31181300,30349100,29633300,29113400,28434800,27238300,26318000,25587900,23277400)
Set-RowValues 77 @(0,0,0,0,0,0,0,0,0,0)
Set-RowValues 80 @(43465,43373,43281,43190,43100,43008,42916,42825,42735,42643)
Set-RowValues 81 @(-997700,793300,951400,943800,885600,773300,678600,354800,689800,506000)
Set-RowValues 83 @(826500,806400,815200,836100,833100,836900,788400,836200,793900,1019900)
Set-RowValues 84 @(0,0,0,0,0,0,0,0,0,0)
Set-RowValues 85 @(0,0,0,0,0,0,0,0,0,0)
Set-RowValues 86 @(0,0,0,0,0,0,0,0,0,0)
Set-RowValues 87 @(0,0,0,0,0,0,0,0,0,0)
Set-RowValues 88 @(0,0,0,0,0,0,0,0,0,0)
Set-RowValues 89 @(-87200,701900,1214600,2295800,1760400,1369900,1147500,2221100,835800,1773600)
Set-RowValues 91 @(-860700,-859200,-929600,-5500,-8100,-4000,-11900,-1277100,-889100,-2139000)
Set-RowValues 92 @(0,0,0,0,0,0,0,0,0,0)
Set-RowValues 93 @(0,0,0,0,0,0,0,0,0,0)
Set-RowValues 94 @(-298500,-348300,-939600,-1447200,-721200,-1173700,-945800,-1217900,141000,-1198700)
Set-RowValues 96 @(-349600,700,-349500,400,-305900,600,-305900,400,-257400,400)
Set-RowValues 97 @(0,0,0,0,0,0,0,0,0,0)
Set-RowValues 98 @(0,0,0,0,0,0,0,0,0,0)
Set-RowValues 99 @(0,0,0,0,0,0,0,0,0,0)
Set-RowValues 100 @(-507300,-251400,793500,-1618500,-748700,-709100,170400,-864900,-2288400,989900)
Set-RowValues 101 @(-163400,-900,-15900,-254800,99500,107600,59800,-122100,530100,-69300)
Set-RowValues 102 @(-1056400,101200,1052600,-1024800,390000,-405400,431900,16200,-781600,1495500)
